$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at B to make room for MgmtNIC (shifts C..N -> D..O)
$ws.Columns.Item(2).Insert()

# Column widths (A unaffected; B = new MgmtNIC; J = DriveHW grew wide; O unaffected-width)
$ws.Columns.Item(2).ColumnWidth = 8.6
$ws.Columns.Item(10).ColumnWidth = 73.1

# Header row (row 1)
$ws.Cells.Item(1,1).Value = "ESXHostname"
$ws.Cells.Item(1,2).Value = "MgmtNIC"
$ws.Cells.Item(1,3).Value = "HostIP"
$ws.Cells.Item(1,4).Value = "HostSubnet"
$ws.Cells.Item(1,5).Value = "HostGW"
$ws.Cells.Item(1,6).Value = "HostMgmtVLAN"
$ws.Cells.Item(1,7).Value = "HostDNS1"
$ws.Cells.Item(1,8).Value = "HostDNS2"
$ws.Cells.Item(1,9).Value = "HostPW"
$ws.Cells.Item(1,10).Value = "DriveHW"
$ws.Cells.Item(1,11).Value = "HostDomain"
$ws.Cells.Item(1,12).Value = "LocalUser"
$ws.Cells.Item(1,13).Value = "LocalPW"
$ws.Cells.Item(1,14).Value = "VCSAIPAddr"
$ws.Cells.Item(1,15).Value = "ListOfPhysicalDrives"

# Row 2 (existing row - update values)
$ws.Cells.Item(2,1).Value = "esx01.tataoui.com"
$ws.Cells.Item(2,2).Value = "vmnic2"
$ws.Cells.Item(2,3).Value = "192.168.10.11"
$ws.Cells.Item(2,4).Value = "255.255.255.0"
$ws.Cells.Item(2,5).Value = "192.168.10.2"
$ws.Cells.Item(2,6).Value = 10
$ws.Cells.Item(2,7).Value = "8.8.8.8"
$ws.Cells.Item(2,8).Value = "192.168.30.2"
$ws.Cells.Item(2,9).Value = "VMware123!"
$ws.Cells.Item(2,10).Value = "naa.500a07510c12161b"
$ws.Cells.Item(2,11).Value = "tataoui.com"
$ws.Cells.Item(2,12).Value = "dwcadmin"
$ws.Cells.Item(2,13).Value = "VMware123!"
$ws.Cells.Item(2,14).Value = "192.168.10.40"
$ws.Cells.Item(2,15).Value = "naa.61866d-Repository MZXLR3T2HBLS2D000H3-VCF_1"

# Row 3 (existing row - update values)
$ws.Cells.Item(3,1).Value = "esx02.tataoui.com"
$ws.Cells.Item(3,2).Value = "vmnic2"
$ws.Cells.Item(3,3).Value = "192.168.10.12"
$ws.Cells.Item(3,4).Value = "255.255.255.0"
$ws.Cells.Item(3,5).Value = "192.168.10.2"
$ws.Cells.Item(3,6).Value = 10
$ws.Cells.Item(3,7).Value = "8.8.8.8"
$ws.Cells.Item(3,8).Value = "192.168.30.2"
$ws.Cells.Item(3,9).Value = "VMware123!"
$ws.Cells.Item(3,10).Value = "naa.500a0751095f1055"
$ws.Cells.Item(3,11).Value = "tataoui.com"
$ws.Cells.Item(3,12).Value = "dwcadmin"
$ws.Cells.Item(3,13).Value = "VMware123!"
$ws.Cells.Item(3,14).Value = "192.168.10.40"
$ws.Cells.Item(3,15).Value = "naa.61866d-Repository MZXLR3T2HBLS2D000H3-VCF_2"

# Row 4 (existing row - update values)
$ws.Cells.Item(4,1).Value = "esx03.tataoui.com"
$ws.Cells.Item(4,2).Value = "vmnic4"
$ws.Cells.Item(4,3).Value = "192.168.10.13"
$ws.Cells.Item(4,4).Value = "255.255.255.0"
$ws.Cells.Item(4,5).Value = "192.168.10.2"
$ws.Cells.Item(4,6).Value = 10
$ws.Cells.Item(4,7).Value = "8.8.8.8"
$ws.Cells.Item(4,8).Value = "192.168.30.2"
$ws.Cells.Item(4,9).Value = "VMware123!"
$ws.Cells.Item(4,10).Value = "t10.ATA_____SATA_SSD________________________________21110524000180______"
$ws.Cells.Item(4,11).Value = "tataoui.com"
$ws.Cells.Item(4,12).Value = "dwcadmin"
$ws.Cells.Item(4,13).Value = "VMware123!"
$ws.Cells.Item(4,14).Value = "192.168.10.40"
$ws.Cells.Item(4,15).Value = "SATA_SSD-Repository Samsung_SSD_980-VCF_3"

# Row 5 (existing row - update values)
$ws.Cells.Item(5,1).Value = "esx04.tataoui.com"
$ws.Cells.Item(5,2).Value = "vmnic0"
$ws.Cells.Item(5,3).Value = "192.168.10.14"
$ws.Cells.Item(5,4).Value = "255.255.255.0"
$ws.Cells.Item(5,5).Value = "192.168.10.2"
$ws.Cells.Item(5,6).Value = 10
$ws.Cells.Item(5,7).Value = "8.8.8.8"
$ws.Cells.Item(5,8).Value = "192.168.30.2"
$ws.Cells.Item(5,9).Value = "VMware123!"
$ws.Cells.Item(5,10).Value = "t10.ATA_____SATA_SSD________________________________19080124004062______"
$ws.Cells.Item(5,11).Value = "tataoui.com"
$ws.Cells.Item(5,12).Value = "dwcadmin"
$ws.Cells.Item(5,13).Value = "VMware123!"
$ws.Cells.Item(5,14).Value = "192.168.10.40"
$ws.Cells.Item(5,15).Value = "SATA-Repository Samsung-SSD_VSAN HITACHI-HDD_VSAN"

# Row 6 (existing row - update values)
$ws.Cells.Item(6,1).Value = "esx05.tataoui.com"
$ws.Cells.Item(6,2).Value = "vmnic0"
$ws.Cells.Item(6,3).Value = "192.168.10.15"
$ws.Cells.Item(6,4).Value = "255.255.255.0"
$ws.Cells.Item(6,5).Value = "192.168.10.2"
$ws.Cells.Item(6,6).Value = 10
$ws.Cells.Item(6,7).Value = "8.8.8.8"
$ws.Cells.Item(6,8).Value = "192.168.30.2"
$ws.Cells.Item(6,9).Value = "VMware123!"
$ws.Cells.Item(6,10).Value = "t10.ATA_____SATA_SSD________________________________18082224001121______"
$ws.Cells.Item(6,11).Value = "tataoui.com"
$ws.Cells.Item(6,12).Value = "dwcadmin"
$ws.Cells.Item(6,13).Value = "VMware123!"
$ws.Cells.Item(6,14).Value = "192.168.10.40"
$ws.Cells.Item(6,15).Value = "SATA_SSD-SSD_VM Samsung-SSD_VSAN HITACHI-HDD_VSAN"

# Row 7 (new row - copy formatting from row 6, then set values)
$ws.Range("A6:O6").Copy()
$ws.Range("A7:O7").PasteSpecial(-4122)
$ws.Cells.Item(7,1).Value = "esx06.tataoui.com"
$ws.Cells.Item(7,2).Value = "vmnic0"
$ws.Cells.Item(7,3).Value = "192.168.10.16"
$ws.Cells.Item(7,4).Value = "255.255.255.0"
$ws.Cells.Item(7,5).Value = "192.168.10.2"
$ws.Cells.Item(7,6).Value = 10
$ws.Cells.Item(7,7).Value = "8.8.8.8"
$ws.Cells.Item(7,8).Value = "192.168.30.2"
$ws.Cells.Item(7,9).Value = "VMware123!"
$ws.Cells.Item(7,10).Value = "t10.ATA_____SATA_SSD________________________________20071324000117______"
$ws.Cells.Item(7,11).Value = "tataoui.com"
$ws.Cells.Item(7,12).Value = "dwcadmin"
$ws.Cells.Item(7,13).Value = "VMware123!"
$ws.Cells.Item(7,14).Value = "192.168.10.40"
$ws.Cells.Item(7,15).Value = "SATA_SSD-SSD_VM Samsung-SSD_VSAN HITACHI-HDD_VSAN"

# Row 8 (new row - copy formatting from row 6, then set values)
$ws.Range("A6:O6").Copy()
$ws.Range("A8:O8").PasteSpecial(-4122)
$ws.Cells.Item(8,1).Value = "esx11.tataoui.com"
$ws.Cells.Item(8,2).Value = "vmnic0"
$ws.Cells.Item(8,3).Value = "192.168.10.17"
$ws.Cells.Item(8,4).Value = "255.255.255.0"
$ws.Cells.Item(8,5).Value = "192.168.10.2"
$ws.Cells.Item(8,6).Value = 10
$ws.Cells.Item(8,7).Value = "8.8.8.8"
$ws.Cells.Item(8,8).Value = "192.168.30.2"
$ws.Cells.Item(8,9).Value = "VMware123!"
$ws.Cells.Item(8,10).Value = "mpx.vmhba32:C0:T0:L0"
$ws.Cells.Item(8,11).Value = "tataoui.com"
$ws.Cells.Item(8,12).Value = "dwcadmin"
$ws.Cells.Item(8,13).Value = "VMware123!"
$ws.Cells.Item(8,14).Value = "192.168.10.40"

# Row 9 (new row - copy formatting from row 6, then set values)
$ws.Range("A6:O6").Copy()
$ws.Range("A9:O9").PasteSpecial(-4122)
$ws.Cells.Item(9,1).Value = "esx12.tataoui.com"
$ws.Cells.Item(9,2).Value = "vmnic0"
$ws.Cells.Item(9,3).Value = "192.168.10.18"
$ws.Cells.Item(9,4).Value = "255.255.255.0"
$ws.Cells.Item(9,5).Value = "192.168.10.2"
$ws.Cells.Item(9,6).Value = 10
$ws.Cells.Item(9,7).Value = "8.8.8.8"
$ws.Cells.Item(9,8).Value = "192.168.30.2"
$ws.Cells.Item(9,9).Value = "VMware123!"
$ws.Cells.Item(9,10).Value = "mpx.vmhba32:C0:T0:L0"
$ws.Cells.Item(9,11).Value = "tataoui.com"
$ws.Cells.Item(9,12).Value = "dwcadmin"
$ws.Cells.Item(9,13).Value = "VMware123!"
$ws.Cells.Item(9,14).Value = "192.168.10.40"

# Row 10 (new row - copy formatting from row 6, then set values)
$ws.Range("A6:O6").Copy()
$ws.Range("A10:O10").PasteSpecial(-4122)
$ws.Cells.Item(10,1).Value = "esx13.tataoui.com"
$ws.Cells.Item(10,2).Value = "vmnic0"
$ws.Cells.Item(10,3).Value = "192.168.10.19"
$ws.Cells.Item(10,4).Value = "255.255.255.0"
$ws.Cells.Item(10,5).Value = "192.168.10.2"
$ws.Cells.Item(10,6).Value = 10
$ws.Cells.Item(10,7).Value = "8.8.8.8"
$ws.Cells.Item(10,8).Value = "192.168.30.2"
$ws.Cells.Item(10,9).Value = "VMware123!"
$ws.Cells.Item(10,10).Value = "mpx.vmhba32:C0:T0:L0"
$ws.Cells.Item(10,11).Value = "tataoui.com"
$ws.Cells.Item(10,12).Value = "dwcadmin"
$ws.Cells.Item(10,13).Value = "VMware123!"
$ws.Cells.Item(10,14).Value = "192.168.10.40"
$ws.Cells.Item(10,15).Value = "SATA_SSD-Repository Samsung_SSD-VCF_3"

$excel.CutCopyMode = $false

# Selection per diff (activeCell C12 on the sheet, matching saved view state)
$ws.Range("C12").Select()

